# Auto-generated Excel COM-interop script
# Updates market-data derived columns (H,I,J,K,L,M,N) on several sheets
# to reflect a scheduled market-price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2625.6182
$ws.Range("J17").Value = 2625.6182
$ws.Range("L17").Value = 7876.8546
$ws.Range("N17").Value = -8212.854599999999
$ws.Range("H107").Value = 8706.75
$ws.Range("I107").Value = 10413.5
$ws.Range("J107").Value = 173
$ws.Range("K107").Value = 10413.5
$ws.Range("L107").Value = 173
$ws.Range("M107").Value = -8493.5
$ws.Range("N107").Value = -4013
$ws.Range("H132").Value = 14250.613
$ws.Range("I132").Value = 2335.1804
$ws.Range("J132").Value = 66167.86
$ws.Range("K132").Value = 7005.541200000001
$ws.Range("L132").Value = 198503.58
$ws.Range("M132").Value = -4475.541200000001
$ws.Range("N132").Value = -203563.58
$ws.Range("H137").Value = 3730.84
$ws.Range("I137").Value = 1171.037
$ws.Range("J137").Value = 6735.826
$ws.Range("K137").Value = 3513.111
$ws.Range("L137").Value = 20207.478
$ws.Range("M137").Value = -963.1109999999999
$ws.Range("N137").Value = -25307.478
$ws.Range("H138").Value = 1907.979
$ws.Range("I138").Value = 1034.7407
$ws.Range("J138").Value = 3058.0977
$ws.Range("K138").Value = 3104.2221
$ws.Range("L138").Value = 9174.293099999999
$ws.Range("M138").Value = 2035.7779
$ws.Range("N138").Value = -19454.2931

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9606.437
$ws.Range("I32").Value = 9168.066000000001
$ws.Range("K32").Value = 9168.066000000001
$ws.Range("M32").Value = -8881.066000000001
$ws.Range("H45").Value = 1912
$ws.Range("I45").Value = 1779.6
$ws.Range("J45").Value = 2243
$ws.Range("K45").Value = 1779.6
$ws.Range("L45").Value = 2243
$ws.Range("M45").Value = -1402.6
$ws.Range("N45").Value = -2997
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992
$ws.Range("H95").Value = 38559
$ws.Range("J95").Value = 38559
$ws.Range("L95").Value = 38559
$ws.Range("N95").Value = -44051
$ws.Range("H96").Value = 32414.666
$ws.Range("J96").Value = 32414.666
$ws.Range("L96").Value = 32414.666
$ws.Range("N96").Value = -37906.666
$ws.Range("H132").Value = 12501858
$ws.Range("I132").Value = 19231988
$ws.Range("K132").Value = 57695964
$ws.Range("M132").Value = -57693434

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1695.975
$ws.Range("I107").Value = 1572.1923
$ws.Range("J107").Value = 1925.8572
$ws.Range("K107").Value = 1572.1923
$ws.Range("L107").Value = 1925.8572
$ws.Range("M107").Value = 347.8077000000001
$ws.Range("N107").Value = -5765.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 942.75
$ws.Range("J16").Value = 860.6
$ws.Range("L16").Value = 860.6
$ws.Range("N16").Value = -1434.6
$ws.Range("H31").Value = 2815.87
$ws.Range("I31").Value = 1178.9259
$ws.Range("J31").Value = 3421.3152
$ws.Range("K31").Value = 1178.9259
$ws.Range("L31").Value = 3421.3152
$ws.Range("M31").Value = -883.9259
$ws.Range("N31").Value = -4011.3152
$ws.Range("H34").Value = 2815.87
$ws.Range("I34").Value = 1178.9259
$ws.Range("J34").Value = 3421.3152
$ws.Range("K34").Value = 1178.9259
$ws.Range("L34").Value = 3421.3152
$ws.Range("M34").Value = -976.9259
$ws.Range("N34").Value = -3825.3152
$ws.Range("H99").Value = 4009.5
$ws.Range("J99").Value = 7014
$ws.Range("L99").Value = 7014
$ws.Range("N99").Value = -10010
$ws.Range("H107").Value = 564.087
$ws.Range("I107").Value = 513.3333
$ws.Range("J107").Value = 659.25
$ws.Range("K107").Value = 513.3333
$ws.Range("L107").Value = 659.25
$ws.Range("M107").Value = 1406.6667
$ws.Range("N107").Value = -4499.25
$ws.Range("H113").Value = 942.75
$ws.Range("J113").Value = 860.6
$ws.Range("L113").Value = 860.6
$ws.Range("N113").Value = -5200.6
$ws.Range("H126").Value = 4009.5
$ws.Range("J126").Value = 7014
$ws.Range("L126").Value = 21042
$ws.Range("N126").Value = -25982
$ws.Range("H132").Value = 30541.611
$ws.Range("I132").Value = 1316.3143
$ws.Range("J132").Value = 103604.86
$ws.Range("K132").Value = 3948.9429
$ws.Range("L132").Value = 310814.58
$ws.Range("M132").Value = -1418.9429
$ws.Range("N132").Value = -315874.58
$ws.Range("H134").Value = 484460.06
$ws.Range("I134").Value = 1217.3
$ws.Range("J134").Value = 1558332.9
$ws.Range("K134").Value = 3651.9
$ws.Range("L134").Value = 4674998.699999999
$ws.Range("M134").Value = -1116.9
$ws.Range("N134").Value = -4680068.699999999
$ws.Range("H135").Value = 44248.57
$ws.Range("J135").Value = 44248.57
$ws.Range("L135").Value = 44248.57
$ws.Range("N135").Value = -54388.57

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 83118.16
$ws.Range("I56").Value = 83118.16
$ws.Range("K56").Value = 83118.16
$ws.Range("M56").Value = -82588.16

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2562.516
$ws.Range("I61").Value = 2532.3447
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2532.3447
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -2330.3447
$ws.Range("N61").Value = -3404
$ws.Range("H100").Value = 1416.7059
$ws.Range("I100").Value = 1405.6
$ws.Range("J100").Value = 1500
$ws.Range("K100").Value = 1405.6
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = -864.5999999999999
$ws.Range("N100").Value = -2582
$ws.Range("H113").Value = 2562.516
$ws.Range("I113").Value = 2532.3447
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2532.3447
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -362.3447000000001
$ws.Range("N113").Value = -7340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 92354.86
$ws.Range("J46").Value = 92354.86
$ws.Range("L46").Value = 92354.86
$ws.Range("N46").Value = -92816.86
$ws.Range("H62").Value = 3192.4092
$ws.Range("I62").Value = 2666.5
$ws.Range("J62").Value = 3245
$ws.Range("K62").Value = 2666.5
$ws.Range("L62").Value = 3245
$ws.Range("M62").Value = -2042.5
$ws.Range("N62").Value = -4493
$ws.Range("H65").Value = 3192.4092
$ws.Range("I65").Value = 2666.5
$ws.Range("J65").Value = 3245
$ws.Range("K65").Value = 13332.5
$ws.Range("L65").Value = 16225
$ws.Range("M65").Value = -10212.5
$ws.Range("N65").Value = -22465
$ws.Range("H107").Value = 12501163
$ws.Range("I107").Value = 1260.4
$ws.Range("K107").Value = 3781.2
$ws.Range("M107").Value = -1861.2
$ws.Range("H113").Value = 1300
$ws.Range("I113").Value = 1366.6666
$ws.Range("K113").Value = 4099.9998
$ws.Range("M113").Value = -1929.9998
$ws.Range("H122").Value = 746
$ws.Range("I122").Value = 717.7778
$ws.Range("K122").Value = 2153.3334
$ws.Range("M122").Value = 296.6666
$ws.Range("H123").Value = 42429
$ws.Range("J123").Value = 42429
$ws.Range("L123").Value = 42429
$ws.Range("N123").Value = -52229
$ws.Range("H134").Value = 92354.86
$ws.Range("J134").Value = 92354.86
$ws.Range("L134").Value = 277064.58
$ws.Range("N134").Value = -282134.58
$ws.Range("H138").Value = 38189.1
$ws.Range("J138").Value = 38189.1
$ws.Range("L138").Value = 38189.1

Write-Host "Updated 187 cells across 7 sheets"
